$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("covid19_cases_switzerland")
$ws2 = $wb.Worksheets.Item("covid19_fatalities_switzerland")

# ===================== covid19_fatalities_switzerland (sheet2) =====================
# New row 22 (2020-04-15, date serial 43916)
$ws2.Range("A22").Value = 43916
$ws2.Range("A22").NumberFormat = "yyyy\-mm\-dd;@"

$ws2.Range("B22").Value = 2
$ws2.Range("D22").Value = 2
$ws2.Range("E22").Value = 7
$ws2.Range("F22").Value = 5
$ws2.Range("G22").Value = 8
$ws2.Range("H22").Value = 1
$ws2.Range("I22").Value = 15
$ws2.Range("K22").Value = 6
$ws2.Range("M22").Value = 1
$ws2.Range("N22").Value = 2
$ws2.Range("Q22").Value = 1
$ws2.Range("S22").Value = 1
$ws2.Range("U22").Value = 1
$ws2.Range("V22").Value = 60
$ws2.Range("X22").Value = 21
$ws2.Range("Y22").Value = 14
$ws2.Range("AA22").Value = 7

$ws2.Range("AB22").Formula = "=SUM(B22:AA22)"
$ws2.Range("AB22").NumberFormat = "0"

$ws2.Activate()
$ws2.Range("G28").Select() | Out-Null

# ===================== covid19_cases_switzerland (sheet1) =====================
# Fill in previously-missing AI (C21) and SO (S21) values for row 21
$ws1.Range("C21").Value = 9
$ws1.Range("S21").Value = 141

# AB21 no longer needs to add back S20/C20 now that C21/S21 are populated directly
$ws1.Range("AB21").Formula = "=SUM(B21:AA21)+X20+H20"

# New row 22 (2020-04-15, date serial 43916)
$ws1.Range("A22").Value = 43916
$ws1.Range("A22").NumberFormat = "yyyy\-mm\-dd;@"

$ws1.Range("B22").Value = 323
$ws1.Range("D22").Value = 40
$ws1.Range("E22").Value = 660
$ws1.Range("Q22").Value = 280
$ws1.Range("R22").Value = 35
$ws1.Range("Z22").Value = 87

$ws1.Range("AB22").Formula = "=SUM(B22:AA22)+AA21+Y21+X20+W21+V21+U21+T21+S21+P21+O21+N21+M21+L21+K21+J21+I21+H20+G21+F21+C21"

$ws1.Activate()
$ws1.Range("N25").Select() | Out-Null
